$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.122.53'
$ws.Range('E2').Value = '  -3.55%  '
$ws.Range('D3').Value = '1.601.71'
$ws.Range('E3').Value = '  -2.89%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '1.001'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = '301.09'
$ws.Range('D7').Value = '0.3778'
$ws.Range('E7').Value = '  -3.19%  '
$ws.Range('D8').Value = '0.3651'
$ws.Range('E8').Value = '  -4.15%  '
$ws.Range('D9').Value = '50.08'
$ws.Range('E9').Value = '  -3.84%  '
$ws.Range('D10').Value = '1.256'
$ws.Range('E10').Value = '  -6.62%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = '0.08132'
$ws.Range('E12').Value = '  -3.71%  '
$ws.Range('D13').Value = '23.02'
$ws.Range('E13').Value = '  -3.37%  '
$ws.Range('E14').Value = '  -6.43%  '
$ws.Range('D15').Value = '7.417'
$ws.Range('E15').Value = '  -7.34%  '
$ws.Range('D16').Value = '0.00001255'
$ws.Range('E16').Value = '  -4.07%  '
$ws.Range('D17').Value = '1.604.48'
$ws.Range('E17').Value = '  -2.85%  '
$ws.Range('D18').Value = '91.57'
$ws.Range('E18').Value = '  -2.94%  '
$ws.Range('D19').Value = '0.06850'
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('E20').Value = '  -7.07%  '
$ws.Range('D21').Value = '6.562'
$ws.Range('E21').Value = '  -5.89%  '
$ws.Range('E22').Value = '  -6.96%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  -5.60%  '
$ws.Range('D25').Value = '23.153.09'
$ws.Range('E25').Value = '  -3.38%  '
$ws.Range('D26').Value = '2.338'
$ws.Range('E26').Value = '  -4.26%  '
$ws.Range('D27').Value = '2.727'
$ws.Range('E27').Value = '  -7.36%  '
$ws.Range('D28').Value = '21.10'
$ws.Range('E28').Value = '  -4.34%  '
$ws.Range('D29').Value = '150.18'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').Value = '5.286'
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('D31').Value = '131.94'
$ws.Range('E31').Value = '  -4.38%  '
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('D33').Value = '6.818'
$ws.Range('E33').Value = '  -13.93%  '
$ws.Range('D34').Value = '1.780.23'
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('D35').Value = '0.07696'
$ws.Range('E35').Value = '  -4.35%  '
$ws.Range('D36').Value = '0.9495'
$ws.Range('E36').Value = '  -6.61%  '
$ws.Range('D37').Value = '0.02747'
$ws.Range('E37').Value = '  -5.92%  '
$ws.Range('D38').Value = '6.262'
$ws.Range('E38').Value = '  -6.99%  '
$ws.Range('D39').Value = '0.2545'
$ws.Range('E39').Value = '  -4.82%  '
$ws.Range('D40').Value = '0.08913'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('D41').Value = '10.05'
$ws.Range('E41').Value = '  -6.37%  '
$ws.Range('D42').Value = '1.386'
$ws.Range('E42').Value = '  -2.44%  '
$ws.Range('D43').Value = '0.7107'
$ws.Range('E43').Value = '  -6.45%  '
$ws.Range('D44').Value = '12.70'
$ws.Range('D45').Value = '15.43'
$ws.Range('E45').Value = '  -4.88%  '
$ws.Range('D46').Value = '0.6627'
$ws.Range('E46').Value = '  -4.86%  '
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = '2.304'
$ws.Range('E48').Value = '  -5.94%  '
$ws.Range('D49').Value = '3.980'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('D50').Value = '132.20'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('D51').Value = '0.07941'
$ws.Range('E51').Value = '  -4.47%  '
